## CS0007 Week 7 -> Week 8 deck: insert the new "Activity 3" slide as the
## 2nd slide (pushing "Random numbers" and "Activity 4" down one spot).

$p = $ppt.ActivePresentation

# Insert a new slide in position 2 using the same "Title and Content"
# layout used by the other content slides in this deck (slideLayout2.xml,
# the 2nd layout registered on the slide master).
$newSlide = $p.Slides.Add(2, 2)

# --- Title -----------------------------------------------------------
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Activity 3"

# --- Body content ------------------------------------------------------
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange

$lines = @(
    "Three main actions need to be done:",
    "Generate rolls of the dice (using Random numbers)",
    "Print out results",
    "Decide whether the user won or not",
    "",
    "",
    "",
    "*Look at your flowchart from activity 2 if you are confused on the steps in the game "
)
$body.Text = [string]::Join("`r", $lines)

# Paragraphs 2-8 are second-level (demoted) bullets.
for ($i = 2; $i -le 8; $i++) {
    $body.Paragraphs($i).IndentLevel = 2
}

# The last line is a plain note with no bullet glyph.
$body.Paragraphs(8).ParagraphFormat.Bullet.Visible = 0
